$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on the target columns so values that look
# numeric (e.g. "0.566", "566.42") are stored as text, matching the
# original inline-string cell type.
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("E2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "63.402.20"
$ws.Range("E2").Value = "  -1.74%  "
$ws.Range("D3").Value = "3.392.79"
$ws.Range("E3").Value = "  -0.74%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "566.42"
$ws.Range("E5").Value = "  -1.09%  "
$ws.Range("D6").Value = "155.84"
$ws.Range("E6").Value = "  -0.75%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "3.396.56"
$ws.Range("E8").Value = "  -0.79%  "
$ws.Range("D9").Value = "0.566"
$ws.Range("E9").Value = "  -8.33%  "
$ws.Range("E10").Value = "  +1.00%  "
$ws.Range("D11").Value = "0.118"
$ws.Range("E11").Value = "  -3.48%  "
$ws.Range("D12").Value = "0.423"
$ws.Range("E12").Value = "  -4.19%  "
$ws.Range("D13").Value = "3.983.02"
$ws.Range("E13").Value = "  -0.68%  "
$ws.Range("E14").Value = "  -0.19%  "
$ws.Range("D15").Value = "26.86"
$ws.Range("E15").Value = "  -3.99%  "
$ws.Range("E16").Value = "  -8.99%  "
$ws.Range("D17").Value = "63.521.01"
$ws.Range("E17").Value = "  -1.64%  "
$ws.Range("D18").Value = "3.403.63"
$ws.Range("E18").Value = "  -0.52%  "
$ws.Range("D19").Value = "6.08"
$ws.Range("E19").Value = "  -4.13%  "
$ws.Range("D20").Value = "13.49"
$ws.Range("E20").Value = "  -3.55%  "
$ws.Range("D21").Value = "383.67"
$ws.Range("E21").Value = "  +2.02%  "
$ws.Range("D22").Value = "7.70"
$ws.Range("E22").Value = "  -3.50%  "
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("D24").Value = "71.00"
$ws.Range("E24").Value = "  -1.98%  "
$ws.Range("D25").Value = "0.514"
$ws.Range("E25").Value = "  -7.56%  "
$ws.Range("E26").Value = "  -3.82%  "
$ws.Range("D27").Value = "9.66"
$ws.Range("E27").Value = "  -5.63%  "
$ws.Range("E28").Value = "  +0.32%  "
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("E30").Value = "  -2.48%  "
$ws.Range("E31").Value = "  -8.10%  "
$ws.Range("D32").Value = "1.97"
$ws.Range("E32").Value = "  -2.48%  "
$ws.Range("D33").Value = "22.77"
$ws.Range("E33").Value = "  -1.56%  "
$ws.Range("E34").Value = "  -4.36%  "
$ws.Range("E35").Value = "  -7.19%  "
$ws.Range("D36").Value = "160.82"
$ws.Range("D37").Value = "0.839"
$ws.Range("E37").Value = "  +9.00%  "
$ws.Range("E38").Value = "  -4.25%  "
$ws.Range("D39").Value = "2.809.87"
$ws.Range("E39").Value = "  -1.26%  "
$ws.Range("D40").Value = "25.81"
$ws.Range("E40").Value = "  -3.33%  "
$ws.Range("D41").Value = "42.83"
$ws.Range("E41").Value = "  -0.15%  "
$ws.Range("D42").Value = "0.0717"
$ws.Range("E42").Value = "  -5.70%  "
$ws.Range("E43").Value = "  -7.05%  "
$ws.Range("D44").Value = "25.55"
$ws.Range("E44").Value = "  -3.59%  "
$ws.Range("D45").Value = "4.34"
$ws.Range("E45").Value = "  -5.99%  "
$ws.Range("E46").Value = "  -3.65%  "
$ws.Range("D47").Value = "325.80"
$ws.Range("E47").Value = "  +1.83%  "
$ws.Range("E48").Value = "  +7.55%  "
$ws.Range("D49").Value = "1.02"
$ws.Range("E49").Value = "  -5.20%  "
$ws.Range("E50").Value = "  -5.84%  "
$ws.Range("D51").Value = "6.26"
$ws.Range("E51").Value = "  -4.88%  "

# Reset the style back to the workbook default so we do not leave a
# stray explicit cell style behind (cells originally had none).
$ws.Range("D2:D51").Style = "Normal"
$ws.Range("E2:E51").Style = "Normal"
